# "Changes made after review"
#
# ShareSkill sheet:
#   - L2 (row 2, "Tags"/Credit column) is retagged from "Credit" to a new
#     "Skill-exchange" shared string.
#   - The view's active selection moves from R5 to L5 (and the visible
#     top-left scrolls from column K to column F, which the engine tracks
#     together with the selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ShareSkill")
$ws.Activate()

# Update the cell value; a new shared string ("Skill-exchange") is created
# automatically since no existing entry matches.
$ws.Range("L2").Value = "Skill-exchange"

# Move the selection/viewport to L5 (previously R5).
$ws.Range("L5").Select()
